$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.845.83'
$ws.Range("E2").Value = '  -2.59%  '

$ws.Range("D3").Value = '3.484.57'
$ws.Range("E3").Value = '  -2.15%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.32'
$ws.Range("E5").Value = '  -2.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.82'
$ws.Range("E6").Value = '  -4.88%  '

$ws.Range("D7").Value = '3.483.27'
$ws.Range("E7").Value = '  -2.32%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("E9").Value = '  -1.57%  '

$ws.Range("E10").Value = '  -2.88%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.59'
$ws.Range("E11").Value = '  +3.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.426'
$ws.Range("E12").Value = '  -3.14%  '

$ws.Range("E13").Value = '  -4.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.71'
$ws.Range("E14").Value = '  -4.91%  '

$ws.Range("D15").Value = '4.071.06'
$ws.Range("E15").Value = '  -2.68%  '

$ws.Range("D16").Value = '3.486.78'
$ws.Range("E16").Value = '  -2.64%  '

$ws.Range("D17").Value = '66.878.55'
$ws.Range("E17").Value = '  -3.47%  '

$ws.Range("E18").Value = '  -0.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.49'
$ws.Range("E19").Value = '  -4.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.45'
$ws.Range("E20").Value = '  -3.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.15'
$ws.Range("E21").Value = '  +0.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '440.68'
$ws.Range("E22").Value = '  -3.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.613'
$ws.Range("E23").Value = '  -4.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.56'
$ws.Range("E24").Value = '  +0.84%  '

$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("D26").Value = '3.618.14'
$ws.Range("E26").Value = '  -2.75%  '

$ws.Range("E27").Value = '  -8.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.79'
$ws.Range("E28").Value = '  -7.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.42'
$ws.Range("E29").Value = '  -7.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.50'
$ws.Range("E30").Value = '  -2.77%  '

$ws.Range("E31").Value = '  -5.51%  '

$ws.Range("E32").Value = '  -0.69%  '

$ws.Range("E33").Value = '  +0.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.50'
$ws.Range("E34").Value = '  -2.81%  '

$ws.Range("E35").Value = '  -6.14%  '

$ws.Range("D36").Value = '3.471.90'
$ws.Range("E36").Value = '  -2.80%  '

$ws.Range("E37").Value = '  -6.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.96'
$ws.Range("E38").Value = '  -4.26%  '

$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '176.20'
$ws.Range("E41").Value = '  -1.37%  '

$ws.Range("E42").Value = '  -2.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.14'
$ws.Range("E43").Value = '  -10.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.44'
$ws.Range("E44").Value = '  -2.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.891'
$ws.Range("E45").Value = '  -0.87%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.15'
$ws.Range("E46").Value = '  -5.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.30'
$ws.Range("E47").Value = '  +1.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.24'
$ws.Range("E48").Value = '  -7.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.50'
$ws.Range("E49").Value = '  -4.02%  '

$ws.Range("E50").Value = '  -8.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.988'
$ws.Range("E51").Value = '  -3.72%  '

